# Weekly driver report update for 2025-04-20
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Bad Drivers" table
# Row 3: Intel(R) Wi-Fi 6E AX211 160MHz - 22.190.0.4
$ws.Range("C3").Value = 2802
$ws.Range("D3").Value = 94.7

# Row 4: Intel(R) Wi-Fi 6E AX211 160MHz - 23.40.0.4
$ws.Range("C4").Value = 493
$ws.Range("D4").Value = 96.6

# Row 5: Totals
$ws.Range("C5").Value = 3295

# "Good Drivers" table
# Row 13: Intel(R) Wi-Fi 6E AX211 160MHz - 22.150.3.1 - fill in the
# previously-blank Driver Vintage date. Format the cell as text first so
# the date-like string isn't auto-converted into a date serial number.
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "2022-08-29"
